$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values edited ("Modified Reg iProctor TC's")
$ws.Range("A2").Value = "NLHAe773"
$ws.Range("B2").Value = 23092282
$ws.Range("C2").Value = "kqjahxd94"
$ws.Range("D2").Value = "S!sw3V&2"
$ws.Range("F2").Value = "utJBzSow"
$ws.Range("G2").Value = "cggC"
